$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet right before "总计"
# ---------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($total)
$q1.Name = "2022-Q1"

# Copy header-row (B1:H1) formatting + the bordered/bold column-A style
# from an existing, identically laid-out sheet ("2021-Q4") so the new
# sheet's styling matches the rest of the workbook.
$tmpl = $wb.Worksheets.Item("2021-Q4")
$tmpl.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$tmpl.Range("A2").Copy()
$q1.Range("A2:A3").PasteSpecial(-4122)

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data rows — D,E,F,G hold text-formatted numbers in the source data
# (keeps leading/trailing zeros exactly as reported), so force the
# "@" text number format before writing them.
$q1.Range("D2:G3").NumberFormat = "@"

$q1.Cells.Item(2,1).Value = 0
$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "014273"
$q1.Range("C2").Value = "广发北交所精选两年定开混合A"
$q1.Range("D2").Value = "4.55"
$q1.Range("E2").Value = "52.69"
$q1.Range("F2").Value = "4.13"
$q1.Range("G2").Value = "0.1879"
$q1.Cells.Item(2,8).Value = 7

$q1.Cells.Item(3,1).Value = 1
$q1.Range("B3").NumberFormat = "@"
$q1.Range("B3").Value = "014274"
$q1.Range("C3").Value = "广发北交所精选两年定开混合C"
$q1.Range("D3").Value = "0.92"
$q1.Range("E3").Value = "52.69"
$q1.Range("F3").Value = "4.13"
$q1.Range("G3").Value = "0.0380"
$q1.Cells.Item(3,8).Value = 7

# ---------------------------------------------------------------
# 2) Add the "2022-Q1" summary row at the top of the "总计" sheet's
#    data, pushing the existing rows down by one.
# ---------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

for ($r = 5; $r -ge 2; $r--) {
    $nr = $r + 1
    $tot.Cells.Item($nr, 2).Value = $tot.Cells.Item($r, 2).Value()
    $tot.Cells.Item($nr, 3).Value = $tot.Cells.Item($r, 3).Value()
    $tot.Cells.Item($nr, 4).Value = $tot.Cells.Item($r, 4).Value()
}

# Re-apply the bordered/bold column-A style down through the newly
# added row 6, then (re)number the index column 0..4.
$tot.Range("A2").Copy()
$tot.Range("A2:A6").PasteSpecial(-4122)

$tot.Cells.Item(2,1).Value = 0
$tot.Cells.Item(3,1).Value = 1
$tot.Cells.Item(4,1).Value = 2
$tot.Cells.Item(5,1).Value = 3
$tot.Cells.Item(6,1).Value = 4

$tot.Cells.Item(2,2).Value = "2022-Q1"
$tot.Cells.Item(2,3).Value = 2
$tot.Cells.Item(2,4).Value = 0.23
